$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Match the formatting of an already-completed row (row 18) for the
# Listening/Reading score cells being filled in on row 19.
$ws.Range("F18").Copy()
$ws.Range("F19").PasteSpecial(-4122)
$ws.Range("H18").Copy()
$ws.Range("H19").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Record the new practice-test results (Cambridge 14, Test 3) on row 19.
$ws.Range("E19").Value = 29
$ws.Range("F19").Formula = '=IFERROR(INDEX(Sheet2!$F$5:$F$20, MATCH(Table1[[#This Row],[Lis_Mark]], Sheet2!$D$5:$D$20, 1)),"No Grade")'
$ws.Range("G19").Value = 25
$ws.Range("H19").Formula = '=IFERROR(INDEX(Sheet2!$F$5:$F$20, MATCH(Table1[[#This Row],[Read_Mark]], Sheet2!$D$5:$D$20, 1)),"No Grade")'
$ws.Range("K19").Formula = '=(F19+H19+I19+J19)/4'

# Clear the placeholder date sitting in the next (still-empty) row.
$ws.Range("C20").ClearContents()

# Leave the selection where the user's cursor ended up after the edit.
$ws.Range("C20").Select()
